$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new column before column B to make room for "Apellido",
# shifting the existing "Partido" / "NIF" columns (and their data) right.
$ws.Columns("B").Insert()

# Rename the first header from "Candidato" to "Nombre" and add the new
# "Apellido" header (write the new "Apellido" header first so it lands
# before "Nombre" in the shared-string table, matching the saved file).
$ws.Range("B1").Value = "Apellido"
$ws.Range("A1").Value = "Nombre"

# Fill in the new surname column alongside the existing first names.
$ws.Range("B3").Value = "Alvarez"
$ws.Range("B4").Value = "Baston"
$ws.Range("B5").Value = "Cienfuegos"

# Update the active selection to match the saved view state.
$ws.Range("E6").Select()
